$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JatMagus")

$ws.Range("A34").Value = "'krjthznfnf"
$ws.Range("B34").Value = "'Mágus"
$ws.Range("C34").Value = "'éjtalizmán"
$ws.Range("D34").Value = "'/Images/Karakterek/magus0.png"
